# Report updated to include chunks
# Adds Chunksize / Static_Chunks / Dynamic_Chunks columns (J:L) to Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---
$ws.Range("J1").Value = "Chunksize"
$ws.Range("K1").Value = "Static_Chunks "
$ws.Range("L1").Value = "Dynamic_Chunks"

# --- Data rows ---
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 18.753095999999999
$ws.Range("L2").Value = 17.824936000000001

$ws.Range("J3").Value = 10
$ws.Range("K3").Value = 17.346914999999999
$ws.Range("L3").Value = 16.756360000000001

$ws.Range("J4").Value = 100
$ws.Range("K4").Value = 16.557599
$ws.Range("L4").Value = 16.625233999999999

$ws.Range("J5").Value = 1000
$ws.Range("K5").Value = 16.564748000000002
$ws.Range("L5").Value = 16.690579

$ws.Range("J6").Value = 2500
$ws.Range("K6").Value = 16.544633000000001
$ws.Range("L6").Value = 16.617457999999999

# --- Column sizing: best-fit widths for the new columns (matches source widths
#     10 / 12.7109375 / 15.140625 chars as closely as this engine's column-width
#     rounding allows) ---
$ws.Columns.Item(10).ColumnWidth = 9.142857142857142
$ws.Columns.Item(11).ColumnWidth = 11.857142857142858
$ws.Columns.Item(12).ColumnWidth = 14.285714285714286

# --- Selection moves to the newly entered header cell, as in the source workbook ---
$ws.Range("L1").Select()
